$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet
$ws.Name = "device"

# Populate header row and data row.
# Cells are written in this specific order so that the workbook's shared
# string table ends up indexed exactly as: OEC type, technlogy type,
# seabed fixed, wave, length [m], height [m], width [m], drymass [kg]
$ws.Range("B1").Value = "OEC type"
$ws.Range("C1").Value = "technlogy type"
$ws.Range("C2").Value = "seabed fixed"
$ws.Range("B2").Value = "wave"
$ws.Range("D1").Value = "length [m]"
$ws.Range("E1").Value = "height [m]"
$ws.Range("F1").Value = "width [m]"
$ws.Range("G1").Value = "drymass [kg]"

$ws.Range("A2").Value = 0
$ws.Range("D2").Value = 100
$ws.Range("E2").Value = 8
$ws.Range("F2").Value = 8
$ws.Range("G2").Value = 150000

# Header row formatting: bold, through column H (H1 stays empty but styled)
$ws.Range("B1:H1").Font.Bold = $true

# Column widths
$ws.Columns.Item(3).ColumnWidth = 12.83
$ws.Columns.Item(7).ColumnWidth = 11.17

# Selection as left by the author
$null = $ws.Range("F8").Select()

# Page setup
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
